# Adapt column header formatting to respective input file names:
#  - "_old" -> "_FV2404"
#  - "_new" -> "_FV2410"
# Freeze the header row, and wrap the data range in an Excel Table (Table1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2404 = "_FV2404"
$fv2410 = "_FV2410"

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J: "<Name>_old" -> "<Name>_FV2404"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + $fv2404
}

# Column K stays "diff"
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U: "<Name>_new" -> "<Name>_FV2410"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + $fv2410
}

# Freeze the header row (row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into an Excel Table named "Table1"
$dataRange = $ws.Range("A1:U61")
$table = $ws.ListObjects.Add(1, $dataRange, $false, 1, $null)
$table.Name = "Table1"
